# Update the lattice multiplication exercises table with new problems.
$d = $word.ActiveDocument
$vt = [char]11  # vertical-tab == Word manual line break (w:br)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "49 x 16" + $vt + "  1    6" + $vt + "  ----" + $vt + "4|    |" + $vt + "9|    |"
$t.Cell(1, 2).Range.Text = "23 x 14" + $vt + "  1    4" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"
$t.Cell(1, 3).Range.Text = "11 x 33" + $vt + "  3    3" + $vt + "  ----" + $vt + "1|    |" + $vt + "1|    |"

$t.Cell(2, 1).Range.Text = "45 x 79" + $vt + "  7    9" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"
$t.Cell(2, 2).Range.Text = "68 x 35" + $vt + "  3    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "8|    |"
$t.Cell(2, 3).Range.Text = "26 x 50" + $vt + "  5    0" + $vt + "  ----" + $vt + "2|    |" + $vt + "6|    |"

$t.Cell(3, 1).Range.Text = "92 x 32" + $vt + "  3    2" + $vt + "  ----" + $vt + "9|    |" + $vt + "2|    |"
$t.Cell(3, 2).Range.Text = "87 x 71" + $vt + "  7    1" + $vt + "  ----" + $vt + "8|    |" + $vt + "7|    |"
$t.Cell(3, 3).Range.Text = "10 x 38" + $vt + "  3    8" + $vt + "  ----" + $vt + "1|    |" + $vt + "0|    |"

$t.Cell(4, 1).Range.Text = "42 x 70" + $vt + "  7    0" + $vt + "  ----" + $vt + "4|    |" + $vt + "2|    |"
$t.Cell(4, 2).Range.Text = "41 x 37" + $vt + "  3    7" + $vt + "  ----" + $vt + "4|    |" + $vt + "1|    |"
$t.Cell(4, 3).Range.Text = "58 x 27" + $vt + "  2    7" + $vt + "  ----" + $vt + "5|    |" + $vt + "8|    |"

$t.Cell(5, 1).Range.Text = "79 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "7|    |" + $vt + "9|    |"
$t.Cell(5, 2).Range.Text = "29 x 13" + $vt + "  1    3" + $vt + "  ----" + $vt + "2|    |" + $vt + "9|    |"
$t.Cell(5, 3).Range.Text = "54 x 25" + $vt + "  2    5" + $vt + "  ----" + $vt + "5|    |" + $vt + "4|    |"

